# Apply changes described by the commit:
#  - Ensure ID of charting date recorded program data element is constant
#    (testchartcode0 -> PatientChartingDate) on the "Test Chart" sheet.
#  - Ensure ID of all complex chart core questions is constant
#    (testchartcorecodeN -> matches the name column) on the "Core" sheet.
#  - Minor formatting touch-ups that tagged along with those cell edits.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$core = $wb.Worksheets.Item("Core")
$testChart = $wb.Worksheets.Item("Test Chart")

# Grab a copy of Metadata!G7's current formatting before we touch it, so we
# can later give it to Test Chart!A2 (which picks up that same look).
$metadata.Range("G7").Copy()
$testChart.Range("A2").PasteSpecial(-4122)

# Normalize Metadata!G7 back to the plain formatting used by the rest of
# the sheet (e.g. matching G8).
$metadata.Range("G8").Copy()
$metadata.Range("G7").PasteSpecial(-4122)

# Core sheet: the "code" column (A) for the complex chart core questions
# now reuses the same constant values as the "name" column (B), instead of
# the old placeholder testchartcorecodeN values.
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartType"
$core.Range("A5").Value = "ComplexChartSubtype"

# Test Chart sheet: the charting date question's code becomes a constant.
$testChart.Range("A2").Value = "PatientChartingDate"
